# Apply "Map 2 and combat system Update 1" changes to the "Phase2" worksheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Phase2")

# New to-do items in column F, rows 26, 28, 29, 30
# (write order matters for shared-string table ordering, so match the
# order in which these distinct strings were first introduced)
$ws.Range("F26").Value = "started implementing chap 2 and combat system"
$ws.Range("F29").Value = "Refactor Player class to not use tilemap"
$ws.Range("F30").Value = "reInstanciate in GameWindow"
$ws.Range("F28").Value = "Player not moving to the right in map 2 because tilemap code is not working on that map"

# Mark rows 35 and 60 ("Status" column D) as Done
$ws.Range("D35").Value = "Done"
$ws.Range("D60").Value = "Done"

# Update the view to reflect scrolled position / new selection
$ws.Activate()
$ws.Range("F31").Select()
$excel.ActiveWindow.ScrollRow = 13
